$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1-General")
Write-Host "Sheet name: $($ws.Name)"
$ws.Range("B3").Value = 2024
